$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11601.8907709685
$ws.Range("D2").Value = 15195.86
$ws.Range("F2").Value = 166.247812560135

$ws.Range("C3").Value = 11162.4351932805
$ws.Range("F3").Value = 305.260151536233

$ws.Range("C4").Value = 7990.17136940428
$ws.Range("F4").Value = 167.625277317654

$ws.Range("C5").Value = 8079.94815841125
$ws.Range("F5").Value = 266.686340005786

$ws.Range("C6").Value = 11464.3205424114
$ws.Range("F6").Value = 399.877289227947

$ws.Range("C7").Value = 10744.8623211183
$ws.Range("F7").Value = 358.142869395615

$ws.Range("C8").Value = 10398.5933548232
$ws.Range("F8").Value = 343.714995799985

$ws.Range("C9").Value = 10653.854149047
$ws.Range("F9").Value = 354.350862225978

$ws.Range("C10").Value = 9869.85747867579
$ws.Range("F10").Value = 321.684334293844

$ws.Range("C11").Value = 7067.13543966631
$ws.Range("F11").Value = 189.537987975942

$ws.Range("C12").Value = 7092.89901295116
$ws.Range("F12").Value = 190.304841073749

$ws.Range("C13").Value = 10505.2988515913
$ws.Range("F13").Value = 328.586893613109

$ws.Range("C14").Value = 10728.4421087691
$ws.Range("F14").Value = 337.884529328852

$ws.Range("C15").Value = 10770.8177647272
$ws.Range("F15").Value = 339.650181660439
